$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the three mailto hyperlinks with new addresses / display text.
# (Range.Hyperlinks.Delete() clears every hyperlink on the sheet in this
# host, so remove them all once and re-add the three we want.)
$ws.Range("C2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:gyv7xyd7@flymail.tk", "", "", "gyv7xyd7@flymail.tk")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:ymhognitf@emlhub.com", "", "", "ymhognitf@emlhub.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:sooseokkim99@gmail.com", "", "", "sooseokkim99@gmail.com")

# Update the "interest" value for the Kann row (D4): nasa -> wrestling
$ws.Range("D4").Value = "wrestling"
